$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full cross-join table (Sending cluster x Target cluster) for Ntf3->Ntrk2,
# replacing the previous 2-row (FAPs/sCs only) table with a 3-cluster
# (ECs/FAPs/sCs) version per Dr Hou's advice.
$rows = @(
  ,@("ECs","Ntf3","Ntrk2","ECs",2,0.6666666666666666,7.310333666666668,21.931001,0.357182641950452,0.357182641950452,2,0.6666666666666666,0.5561916666666666,1.668575,0.01178629317764927,0.01178629317764927,4.065946665952779,36.593519993575,0.004209859335995355,0.004209859335995355)
  ,@("ECs","Ntf3","Ntrk2","FAPs",2,0.6666666666666666,7.310333666666668,21.931001,0.357182641950452,0.357182641950452,3,1,43.12631833333334,129.378955,0.9138925697963186,0.9138925697963186,315.2677768315506,2837.409991483956,0.326426562538737,0.326426562538737)
  ,@("ECs","Ntf3","Ntrk2","sCs",2,0.6666666666666666,7.310333666666668,21.931001,0.357182641950452,0.357182641950452,3,1,3.507192333333334,10.521577,0.0743211370260321,0.07432113702603209,25.6387461898419,230.748715708577,0.02654622007571971,0.0265462200757197)
  ,@("FAPs","Ntf3","Ntrk2","ECs",3,1,10.17625966666667,30.528779,0.4972116839874968,0.4972116839874969,2,0.6666666666666666,0.5561916666666666,1.668575,0.01178629317764927,0.01178629317764927,5.659950824436111,50.93955741992499,0.00586028267882934,0.00586028267882934)
  ,@("FAPs","Ntf3","Ntrk2","FAPs",3,1,10.17625966666667,30.528779,0.4972116839874968,0.4972116839874969,3,1,43.12631833333334,129.378955,0.9138925697963186,0.9138925697963186,438.8646138273273,3949.781524445946,0.4543980636120886,0.4543980636120886)
  ,@("FAPs","Ntf3","Ntrk2","sCs",3,1,10.17625966666667,30.528779,0.4972116839874968,0.4972116839874969,3,1,3.507192333333334,10.521577,0.0743211370260321,0.07432113702603209,35.69009988494256,321.210898964483,0.03695333769657892,0.03695333769657892)
  ,@("sCs","Ntf3","Ntrk2","ECs",2,0.6666666666666666,2.980061000000001,8.940183000000001,0.1456056740620512,0.1456056740620512,2,0.6666666666666666,0.5561916666666666,1.668575,0.01178629317764927,0.01178629317764927,1.657485094358333,14.917365849225,0.001716151162824577,0.001716151162824577)
  ,@("sCs","Ntf3","Ntrk2","FAPs",2,0.6666666666666666,2.980061000000001,8.940183000000001,0.1456056740620512,0.1456056740620512,3,1,43.12631833333334,129.378955,0.9138925697963186,0.9138925697963186,128.5190593387517,1156.671534048765,0.1330679436454931,0.1330679436454931)
  ,@("sCs","Ntf3","Ntrk2","sCs",2,0.6666666666666666,2.980061000000001,8.940183000000001,0.1456056740620512,0.1456056740620512,3,1,3.507192333333334,10.521577,0.0743211370260321,0.07432113702603209,10.45164709206567,94.06482382859102,0.01082157925373347,0.01082157925373347)
)

$r = 2
foreach ($row in $rows) {
  $c = 1
  foreach ($val in $row) {
    $ws.Cells.Item($r, $c).Value = $val
    $c = $c + 1
  }
  $r = $r + 1
}
